$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing meanrank values for rows 2-3 (RF, GNN-MT)
$ws.Range("B2").Value = 3.442857142857143
$ws.Range("B3").Value = 3.216666666666667

# Row 4 now becomes the new "GNN-MT-O" variant (reusing the old PN slot)
$ws.Range("A4").Value = "8_train (GNN-MT-O) val delta-auprc"
$ws.Range("B4").Value = 3.016666666666667

# Row 5: the original "PN" label, now with an updated value, placed after GNN-MT-O
$ws.Range("A5").Value = "8_train (PN) val delta-auprc"
$ws.Range("B5").Value = 2.642857142857143

# Row 6: new "PN-O" variant
$ws.Range("A6").Value = "8_train (PN-O) val delta-auprc"
$ws.Range("B6").Value = 2.285714285714286

# Apply the same label formatting (bold, border, centered) used for the
# other label cells in column A to the two newly added rows.
$ws.Range("A2").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
